$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NCT(2.733499465052115, 1.0124218441036286, 0.2801831469448661, 2.1472377843122805)"
$ws.Range("C2").Value = "NIG(0.6314786235156741, 0.4416890479884843, 5.649031923965321, 4.7667876086525505)"
$ws.Range("D2").Value = "JSB(11.702062825044766, 2.184749358821729, -6.902419207237706, 2265.562738089097)"
$ws.Range("E2").Value = "NIG(1.8167914545491521, 1.4518001816670671, 3.4378365801249053, 6.441663861112954)"
